# Generate Report for Handoff
# Rewrites the localization-status report for the new source file
# (435fd151-6499-403e-946d-f595acc6a1b7.md) and its freshly generated
# handoff xliffs, clearing out the stale handback bookkeeping for the
# two locales that haven't come back yet.

$wb = $excel.ActiveWorkbook

$oldGuid = "c15e0277-08c4-42d0-9f6c-2f5a303b825f"
$newGuid = "435fd151-6499-403e-946d-f595acc6a1b7"
$oldHash = "ad6d05f0e6091d8398d032365659974c82ec8922"
$newHash = "08f0fff93682c32190375061a122d6ac3a93917e"

$newFileName = $newGuid + ".md"
$newPathAndName = "e2e\" + $newGuid + ".md"
$newHoDate = "2016-08-13 03:13:20"

$newZhHandoffFile = $newGuid + "." + $newHash + ".zh-cn.xlf"
$newZhHandoffDate = "2016-08-13 03:13:13"

$newDeHandoffFile = $newGuid + "." + $newHash + ".de-de.xlf"

$clearedHandbackDate = "0001-01-01 00:00:00"

function Remove-CellHyperlink($ws, $addr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.Delete()
        }
    }
}

function Clear-TargetAndHandback($ws) {
    # Latest Target File / Latest Handback File: no handback has come in yet
    # for the new source file, so these go back to "not set" (blank text)
    # and the stale hyperlink + hyperlink styling on the target-file cell
    # is removed.
    Remove-CellHyperlink $ws '$I$2'

    $cellI = $ws.Range("I2")
    $cellI.Font.Underline = -4142
    $cellI.Font.ColorIndex = -4105
    $cellI.Value = "'"

    $ws.Range("J2").Value = "'"
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathAndName
$wsOverview.Hyperlinks.Item(1).TextToDisplay = $newPathAndName
$wsOverview.Range("G2").Value = $newHoDate

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newFileName
$wsZh.Range("G2").Value = $newZhHandoffFile
$wsZh.Range("H2").Value = $newZhHandoffDate
$wsZh.Range("K2").Value = $clearedHandbackDate

Clear-TargetAndHandback $wsZh

$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newFileName
$wsDe.Range("G2").Value = $newDeHandoffFile
$wsDe.Range("H2").Value = $newHoDate
$wsDe.Range("K2").Value = $clearedHandbackDate

Clear-TargetAndHandback $wsDe

$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426
